$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed values from the
# coinranking.com symbol-list sync. Values are written with a leading apostrophe
# so Excel stores/keeps them as literal text (matching the original inline-string
# cells) instead of auto-converting numeric-looking / percent-looking text to a
# number.
$ws.Range("D2").Value = "'302.97"
$ws.Range("E2").Value = "'2.64%"
$ws.Range("D3").Value = "'43.22"
$ws.Range("E3").Value = "'6.45%"
$ws.Range("D4").Value = "'5.065"
$ws.Range("E4").Value = "'0.27%"
$ws.Range("D5").Value = "'0.07686"
$ws.Range("E5").Value = "'3.59%"
$ws.Range("E6").Value = "'3.47%"
$ws.Range("E7").Value = "'12.03%"
$ws.Range("D8").Value = "'0.1260"
$ws.Range("E8").Value = "'4.97%"
$ws.Range("D9").Value = "'0.1855"
$ws.Range("E9").Value = "'2.78%"
$ws.Range("D10").Value = "'0.09165"
$ws.Range("E10").Value = "'4.05%"
$ws.Range("E11").Value = "'-3.55%"
$ws.Range("D12").Value = "'0.1049"
$ws.Range("E12").Value = "'0.08%"
$ws.Range("D13").Value = "'0.001281"
$ws.Range("E13").Value = "'1.40%"
$ws.Range("D14").Value = "'0.005745"
$ws.Range("E14").Value = "'-3.58%"
$ws.Range("E15").Value = "'1,894.53%"
$ws.Range("D16").Value = "'3.344"
$ws.Range("E16").Value = "'-0.45%"
$ws.Range("D17").Value = "'4.413"
$ws.Range("E17").Value = "'1.33%"
$ws.Range("E18").Value = "'-1.96%"
$ws.Range("E19").Value = "'1.35%"
$ws.Range("D20").Value = "'8.651"
$ws.Range("E20").Value = "'8.89%"
$ws.Range("E21").Value = "'-0.66%"
$ws.Range("E22").Value = "'7.66%"
$ws.Range("D23").Value = "'0.04173"
$ws.Range("E23").Value = "'4.12%"
$ws.Range("D24").Value = "'0.001283"
$ws.Range("E24").Value = "'0.91%"
$ws.Range("D25").Value = "'0.004464"
$ws.Range("E25").Value = "'15.29%"
$ws.Range("D26").Value = "'0.0001348"
$ws.Range("E26").Value = "'9.47%"
$ws.Range("D38").Value = "'0.02456"
$ws.Range("E38").Value = "'3.79%"
$ws.Range("D39").Value = "'0.05277"
$ws.Range("E39").Value = "'1.82%"
$ws.Range("D40").Value = "'0.005941"
$ws.Range("E40").Value = "'-0.58%"
$ws.Range("D41").Value = "'0.007670"
$ws.Range("E41").Value = "'-1.32%"
$ws.Range("D42").Value = "'0.1348"
$ws.Range("E42").Value = "'2.19%"
$ws.Range("D43").Value = "'0.007371"
$ws.Range("E43").Value = "'-0.16%"
$ws.Range("D44").Value = "'0.007564"
$ws.Range("E44").Value = "'5.03%"
$ws.Range("D45").Value = "'0.3011"
$ws.Range("E45").Value = "'2.56%"
$ws.Range("D46").Value = "'0.00006696"
$ws.Range("E46").Value = "'5.78%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.28%"
$ws.Range("D48").Value = "'0.03293"
$ws.Range("E48").Value = "'-29.25%"
$ws.Range("E49").Value = "'-0.08%"
$ws.Range("D50").Value = "'0.00002097"
$ws.Range("E50").Value = "'-0.28%"
$ws.Range("D51").Value = "'0.0001997"
$ws.Range("E51").Value = "'-0.28%"
